# Apply odds/value updates for the 2025-05-05 FlashScore weekly fixtures sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.2
$ws.Range("T2").Value = 9
$ws.Range("AD2").Value = 251
$ws.Range("AJ2").Value = 34

# Row 3
$ws.Range("N3").Value = 1.91
$ws.Range("O3").Value = 1.99

# Row 10
$ws.Range("N10").Value = 2.1
$ws.Range("O10").Value = 1.7

# Row 16
$ws.Range("G16").Value = 3.6
$ws.Range("H16").Value = 3.8
$ws.Range("I16").Value = 1.9
$ws.Range("U16").Value = 21
$ws.Range("AE16").Value = 9
$ws.Range("AF16").Value = 10
$ws.Range("AG16").Value = 8.5
$ws.Range("AH16").Value = 17

# Row 17
$ws.Range("G17").Value = 2
$ws.Range("I17").Value = 3.4
$ws.Range("U17").Value = 11
$ws.Range("V17").Value = 9
$ws.Range("W17").Value = 19
$ws.Range("AE17").Value = 12

# Row 18
$ws.Range("N18").Value = 1.44
$ws.Range("O18").Value = 2.7

# Row 19
$ws.Range("K19").Value = 19
$ws.Range("L19").Value = 1.14
$ws.Range("M19").Value = 5.5
$ws.Range("N19").Value = 1.48
$ws.Range("O19").Value = 2.6

# Row 20
$ws.Range("G20").Value = 2.35
$ws.Range("I20").Value = 2.9
$ws.Range("X20").Value = 17

# Row 21
$ws.Range("G21").Value = 3.3
$ws.Range("L21").Value = 1.2
$ws.Range("M21").Value = 4.33
$ws.Range("N21").Value = 1.67
$ws.Range("O21").Value = 2.15
$ws.Range("R21").Value = 1.57
$ws.Range("S21").Value = 2.25

# Row 29
$ws.Range("G29").Value = 3.55
$ws.Range("H29").Value = 3.55
$ws.Range("L29").Value = 1.21
$ws.Range("M29").Value = 3.5
$ws.Range("O29").Value = 2
$ws.Range("W29").Value = 50
$ws.Range("Z29").Value = 12.5
$ws.Range("AA29").Value = 7
$ws.Range("AB29").Value = 13
$ws.Range("AD29").Value = 300
$ws.Range("AH29").Value = 17

# Row 36
$ws.Range("N36").Value = 2.2
$ws.Range("O36").Value = 1.65

# Row 40
$ws.Range("L40").Value = 1.31
$ws.Range("M40").Value = 2.9
$ws.Range("N40").Value = 1.9
$ws.Range("P40").Value = 1.42
$ws.Range("Q40").Value = 2.65
$ws.Range("R40").Value = 1.72
$ws.Range("S40").Value = 1.88
$ws.Range("T40").Value = 8.5
$ws.Range("U40").Value = 13.5
$ws.Range("V40").Value = 10
$ws.Range("X40").Value = 23
$ws.Range("Y40").Value = 32
$ws.Range("Z40").Value = 9.5
$ws.Range("AB40").Value = 14
$ws.Range("AC40").Value = 65
$ws.Range("AD40").Value = 500
$ws.Range("AE40").Value = 7.9
$ws.Range("AI40").Value = 20
$ws.Range("AJ40").Value = 30

# Row 41
$ws.Range("G41").Value = 1.72
$ws.Range("H41").Value = 3.65
$ws.Range("I41").Value = 4.25
$ws.Range("L41").Value = 1.21
$ws.Range("M41").Value = 3.5
$ws.Range("O41").Value = 2
$ws.Range("U41").Value = 9
$ws.Range("V41").Value = 8
$ws.Range("W41").Value = 14.5
$ws.Range("X41").Value = 12.5
$ws.Range("Y41").Value = 21
$ws.Range("AA41").Value = 7.3
$ws.Range("AB41").Value = 14
$ws.Range("AE41").Value = 13.5
$ws.Range("AF41").Value = 26
$ws.Range("AG41").Value = 14
$ws.Range("AH41").Value = 70
$ws.Range("AI41").Value = 40
$ws.Range("AJ41").Value = 40

# Row 42
$ws.Range("J42").Value = 1.07
$ws.Range("K42").Value = 9
$ws.Range("N42").Value = 2.15
$ws.Range("O42").Value = 1.67

# Row 44
$ws.Range("J44").Value = 1.07
$ws.Range("K44").Value = 9
$ws.Range("N44").Value = 2.15
$ws.Range("O44").Value = 1.67

# Row 51
$ws.Range("G51").Value = 2.15
$ws.Range("H51").Value = 3
$ws.Range("I51").Value = 3.45
$ws.Range("M51").Value = 2.12
$ws.Range("N51").Value = 2.62
$ws.Range("O51").Value = 1.38
$ws.Range("P51").Value = 1.6
$ws.Range("Q51").Value = 2.05
$ws.Range("R51").Value = 2.25
$ws.Range("S51").Value = 1.5
$ws.Range("T51").Value = 5.2
$ws.Range("U51").Value = 8.5
$ws.Range("V51").Value = 10
$ws.Range("W51").Value = 20
$ws.Range("X51").Value = 24
$ws.Range("Y51").Value = 55
$ws.Range("Z51").Value = 5.9
$ws.Range("AA51").Value = 6.2
$ws.Range("AB51").Value = 23
$ws.Range("AC51").Value = 175
$ws.Range("AE51").Value = 7
$ws.Range("AF51").Value = 15.5
$ws.Range("AG51").Value = 13.5
$ws.Range("AH51").Value = 50
$ws.Range("AI51").Value = 45
$ws.Range("AJ51").Value = 70

# Row 58
$ws.Range("J58").Value = 1.05
$ws.Range("K58").Value = 11

# Row 65
$ws.Range("G65").Value = 1.75
$ws.Range("H65").Value = 3.7
$ws.Range("I65").Value = 4.33
$ws.Range("N65").Value = 1.95
$ws.Range("O65").Value = 1.9
$ws.Range("W65").Value = 15
$ws.Range("Z65").Value = 10
$ws.Range("AI65").Value = 34

# Row 67
$ws.Range("G67").Value = 1.38
$ws.Range("H67").Value = 5
$ws.Range("L67").Value = 1.18
$ws.Range("M67").Value = 4.5
$ws.Range("N67").Value = 1.62
$ws.Range("O67").Value = 2.25
$ws.Range("AJ67").Value = 51

# Row 68
$ws.Range("G68").Value = 2
$ws.Range("I68").Value = 3.5
$ws.Range("P68").Value = 1.36
$ws.Range("Q68").Value = 3
$ws.Range("W68").Value = 17
$ws.Range("Z68").Value = 11

# Row 71
$ws.Range("G71").Value = 2.3
$ws.Range("I71").Value = 2.63
$ws.Range("L71").Value = 1.14
$ws.Range("M71").Value = 5.5
$ws.Range("R71").Value = 1.44
$ws.Range("S71").Value = 2.63
$ws.Range("Z71").Value = 19
$ws.Range("AG71").Value = 11
$ws.Range("AH71").Value = 29

# Row 73
$ws.Range("G73").Value = 1.93
$ws.Range("H73").Value = 3.7
$ws.Range("I73").Value = 3.4
$ws.Range("P73").Value = 1.34
$ws.Range("Q73").Value = 3
$ws.Range("R73").Value = 1.62
$ws.Range("S73").Value = 2.15
$ws.Range("T73").Value = 8.75
$ws.Range("U73").Value = 10.25
$ws.Range("W73").Value = 17
$ws.Range("X73").Value = 14
$ws.Range("AA73").Value = 7.2
$ws.Range("AB73").Value = 13
$ws.Range("AD73").Value = 350
$ws.Range("AE73").Value = 12.5
$ws.Range("AF73").Value = 19.5
$ws.Range("AG73").Value = 11.75
